$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A2: Pipeline repr containing a NamedFeatureSelector object id -> update memory address
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f7b44610850>),`n                ('model',`n                 LogisticRegression(C=0.001, max_iter=1000, penalty='l1',`n                                    random_state=42, solver='saga'))])"

# Cell C2: Best params dict containing a NamedFeatureSelector object id -> update memory address
$ws.Range("C2").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f7bf754aac0>, 'scaler': MinMaxScaler(), 'model__solver': 'saga', 'model__penalty': 'l1', 'model__class_weight': None, 'model__C': 0.001}"

# Cell A5: Pipeline repr containing a NamedFeatureSelector object id -> update memory address
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', None),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7f7be8494d60>),`n                ('model',`n                 LogisticRegression(C=5, max_iter=1000, penalty='l1',`n                                    random_state=42, solver='saga'))])"

# Cell C5: Best params dict containing a NamedFeatureSelector object id -> update memory address
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7f7bf7ac4a00>, 'scaler': None, 'model__solver': 'saga', 'model__penalty': 'l1', 'model__class_weight': None, 'model__C': 5}"
